# --------------------------------------------------------------------------
# Add a new "2022-Q1" sheet (fund-holder detail) before the "总计" sheet, and
# insert a new summary row at the top of "总计" for the 2022-Q1 period.
# --------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("2021-Q4")
$total    = $wb.Worksheets.Item("总计")

# ---------------------------------------------------------------------
# 1) "总计" sheet: insert a new row 2 and fill it with the 2022-Q1 totals
# ---------------------------------------------------------------------
$total.Rows.Item(2).Insert()

# Re-apply the data-row number format (copied from the row that was pushed
# down to row 3) so the new row matches the existing rows' styling.
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Cells.Item(2, 1).Value = 0
$total.Cells.Item(2, 2).Value = "2022-Q1"
$total.Cells.Item(2, 3).Value = 47
$total.Cells.Item(2, 4).Value = 10.08

# ---------------------------------------------------------------------
# 2) Create the new "2022-Q1" sheet right before "总计"
# ---------------------------------------------------------------------
$newSheet = $wb.Worksheets.Add($total)
$newSheet.Name = "2022-Q1"

# Header row (B1:H1) - copy formatting from the template sheet, then set text
$template.Range("B1:H1").Copy()
$newSheet.Range("B1:H1").PasteSpecial(-4122)

$newSheet.Cells.Item(1, 2).Value = "基金代码"
$newSheet.Cells.Item(1, 3).Value = "基金名称"
$newSheet.Cells.Item(1, 4).Value = "基金规模"
$newSheet.Cells.Item(1, 5).Value = "股票总仓位"
$newSheet.Cells.Item(1, 6).Value = "仓位占比"
$newSheet.Cells.Item(1, 7).Value = "持有市值(亿元)"
$newSheet.Cells.Item(1, 8).Value = "仓位排名"

# Column A (row index) formatting - copy from template's A2 down through A48
$template.Range("A2").Copy()
$newSheet.Range("A2:A48").PasteSpecial(-4122)

# Fund-holder detail rows: code, name, fund size, stock position, position
# ratio, held market value, position rank
$fundData = @(
    @("310358","申万菱信新经济混合","41.92","77.32","4.03","1.6894",7),
    @("506005","博时科创板三年定期开放混合","22.84","96.44","4.18","0.9547",4),
    @("506000","南方科创板 3 年定期开放混合","24.62","96.87","3.58","0.8814",7),
    @("001404","招商移动互联网产业股票","13.45","90.96","5.68","0.7640",3),
    @("004666","长城久嘉创新成长灵活配置混合","26.19","92.18","2.86","0.7490",9),
    @("011488","申万菱信乐享混合","12.18","77.88","4.18","0.5091",6),
    @("013175","海富通碳中和混合A","14.27","93.69","3.09","0.4409",10),
    @("011201","财通优势行业轮动混合A","10.65","88.39","3.95","0.4207",8),
    @("001070","建信信息产业股票","10.67","81.80","3.38","0.3606",8),
    @("012210","申万菱信智能汽车股票型证券投资基金A","4.76","82.52","5.93","0.2823",3),
    @("008655","招商科技创新混合A","4.72","90.71","5.68","0.2681",2),
    @("001150","融通互联网传媒灵活配置混合","9.18","92.66","2.77","0.2543",7),
    @("013176","海富通碳中和混合C","7.99","93.69","3.09","0.2469",10),
    @("012650","博时半导体主题混合型证券投资基金A","7.18","92.40","2.94","0.2111",10),
    @("160919","大成产业升级股票(LOF)","3.95","87.76","3.99","0.1576",10),
    @("519026","海富通中小盘混合","3.88","91.14","4.04","0.1568",8),
    @("008962","建信科技创新混合A","4.34","84.82","3.50","0.1519",9),
    @("012051","申万菱信乐道三年持有期混合型证券投资基金","3.38","81.64","4.02","0.1359",9),
    @("013634","申万菱信双利混合A","7.83","22.26","1.55","0.1214",7),
    @("006281","万家人工智能混合","2.94","90.50","3.93","0.1155",9),
    @("000308","建信创新中国混合","3.11","84.50","3.62","0.1126",7),
    @("012651","博时半导体主题混合型证券投资基金C","3.14","92.40","2.94","0.0923",10),
    @("000522","华润元大信息传媒科技混合","1.50","70.63","5.98","0.0897",5),
    @("001924","华夏国企改革灵活配置混合","2.78","91.54","3.21","0.0892",10),
    @("008656","招商科技创新混合C","1.57","90.71","5.68","0.0892",2),
    @("004314","前海开源沪港深新硬件主题灵活配置混合A","1.67","90.05","5.12","0.0855",7),
    @("012211","申万菱信智能汽车股票型证券投资基金C","1.40","82.52","5.93","0.0830",3),
    @("013339","创金合信芯片产业股票A","1.50","93.43","4.79","0.0718",8),
    @("013345","富荣信息技术混合A","1.96","90.39","3.48","0.0682",8),
    @("673141","西部利得景程灵活配置混合A","2.09","86.50","3.15","0.0658",9),
    @("006502","财通集成电路产业股票A","1.29","79.76","4.62","0.0596",6),
    @("004315","前海开源沪港深新硬件主题灵活配置混合C","1.00","90.05","5.12","0.0512",7),
    @("013346","富荣信息技术混合C","1.44","90.39","3.48","0.0501",8),
    @("501032","财通福盛多策略混合（LOF）","0.74","92.47","4.02","0.0297",10),
    @("013340","创金合信芯片产业股票C","0.61","93.43","4.79","0.0292",8),
    @("006503","财通集成电路产业股票C","0.46","79.76","4.62","0.0213",6),
    @("673143","西部利得景程灵活配置混合C","0.60","86.50","3.15","0.0189",9),
    @("011202","财通优势行业轮动混合C","0.43","88.39","3.95","0.0170",8),
    @("009882","华润元大核心动力混合A","0.22","68.63","6.66","0.0147",1),
    @("006818","安信盈利驱动股票A","0.27","83.93","5.18","0.0140",3),
    @("001574","中海混改红利主题精选灵活配置混合","0.30","89.77","4.07","0.0122",7),
    @("013635","申万菱信双利混合C","0.75","22.26","1.55","0.0116",7),
    @("008963","建信科技创新混合C","0.26","84.82","3.50","0.0091",9),
    @("006819","安信盈利驱动股票C","0.17","83.93","5.18","0.0088",3),
    @("013903","国泰君安信息行业混合","0.25","84.06","3.37","0.0084",6),
    @("009883","华润元大核心动力混合C","0.09","68.63","6.66","0.0060",1),
    @("008890","中邮价值优选一年定期开放灵活配置混合","0.13","62.02","3.08","0.0040",6)
)

$r = 2
foreach ($item in $fundData) {
    $newSheet.Cells.Item($r, 1).Value = $r - 2
    $newSheet.Cells.Item($r, 2).Value = "'" + $item[0]
    $newSheet.Cells.Item($r, 3).Value = $item[1]
    $newSheet.Cells.Item($r, 4).Value = "'" + $item[2]
    $newSheet.Cells.Item($r, 5).Value = "'" + $item[3]
    $newSheet.Cells.Item($r, 6).Value = "'" + $item[4]
    $newSheet.Cells.Item($r, 7).Value = "'" + $item[5]
    $newSheet.Cells.Item($r, 8).Value = $item[6]
    $r = $r + 1
}
